$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "0.475"; New = "0.459" },
    @{ Old = "0.845"; New = "0.887" },
    @{ Old = "0.811"; New = "0.807" },
    @{ Old = "0.534"; New = "0.619" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.New, 2)
}
